$d = $word.ActiveDocument

function Set-ParaXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $xml = "<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">$innerXml</w:p>"
    $p.Range.InsertXML($xml)
}

# Reusable <w:pPr> blocks (list paragraph / numbering properties) used below.
$pPrNum2Ilvl0 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>'
$pPrNum3Ilvl0 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>'
$pPrNum3Ilvl1 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr>'

# --- "Abastract" (paragraph 3) -> wrap run with spell-check proofErr markers ---
Set-ParaXml 3 '<w:proofErr w:type="spellStart"/><w:r><w:t>Abastract</w:t></w:r><w:proofErr w:type="spellEnd"/>'

# --- "Norad" (paragraph 12) -> wrap run with spell-check proofErr markers ---
$noradInner = "${pPrNum2Ilvl0}<w:proofErr w:type=`"spellStart`"/><w:r><w:t>Norad</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>"
Set-ParaXml 12 $noradInner

# --- "Acknolwdments" (paragraph 23) -> wrap run with spell-check proofErr markers ---
$ackInner = "${pPrNum3Ilvl0}<w:proofErr w:type=`"spellStart`"/><w:r><w:t>Acknolwdments</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>"
Set-ParaXml 23 $ackInner

# --- "Cube Sats" (paragraph 27) -> split into "Cube " + proofErr-wrapped "Sats" ---
$cubeSatsInner = "${pPrNum3Ilvl1}<w:r><w:t xml:space=`"preserve`">Cube </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>Sats</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>"
Set-ParaXml 27 $cubeSatsInner

# --- "Desribe how they work and stuff" (paragraph 30) -> split proofErr-wrapped "Desribe" + rest ---
$desribeInner = "${pPrNum3Ilvl1}<w:proofErr w:type=`"spellStart`"/><w:r><w:t>Desribe</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> how they work and stuff</w:t></w:r>"
Set-ParaXml 30 $desribeInner

# --- "Querys " (paragraph 33) -> becomes "Explain they give SATCAT and TLE", ---
# --- and a brand-new paragraph is inserted after it holding the old "Querys " text ---
$explainInner = "${pPrNum3Ilvl1}<w:r><w:t>Explain they give SATCAT and TLE</w:t></w:r>"
Set-ParaXml 33 $explainInner

$p33 = $d.Paragraphs(33)
$p33.Range.InsertParagraphAfter()

$querysInner = "${pPrNum3Ilvl1}<w:proofErr w:type=`"spellStart`"/><w:r><w:t>Querys</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r>"
Set-ParaXml 34 $querysInner

# --- "The Code" (now paragraph 35) -> remove the _GoBack bookmark pair ---
$theCodeInner = "${pPrNum3Ilvl0}<w:r><w:t>The Code</w:t></w:r>"
Set-ParaXml 35 $theCodeInner

# --- "Use ciataiton file on github" (now paragraph 39) -> split into several runs, ---
# --- wrapping the misspelled words with proofErr markers ---
$citationInner = '<w:r><w:t xml:space="preserve">Use </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ciataiton</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> file on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>github</w:t></w:r><w:proofErr w:type="spellEnd"/>'
Set-ParaXml 39 $citationInner

# --- Final (now empty) paragraph (now paragraph 45) -> gets the large new block of ---
# --- text describing Get_SATCAT.m, with the _GoBack bookmark re-added inside it ---
$pPrLast = '<w:pPr><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr></w:pPr>'
$tailFrag = '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>Get_SATCAT.m</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> is the </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>MATLAB</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> file that gets the satellite </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">catalog </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> numbers of all orbital debris </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>launched</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> after a given year and with the </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>“</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>R</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>CS_SIZE” value equal to “SMALL”</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">. They want </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>your</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> can be determined by the user''s input by default is sets 1990 however if you desire more or less information it can be adjusted. No hear that putting the launch here as earlier than 1990 can cause a timeout error as the values to a large if this occurs the timeout value should be adjusted to be longer.</w:t></w:r>'
$lastInner = $pPrLast + $tailFrag
Set-ParaXml 45 $lastInner
